# fdo#75168 - add a second sheet with "expression" type conditional
# formatting rules, exercising cfRule formula load/save round-tripping.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, placed right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Data
$ws2.Range("A1").Value = 2
$ws2.Range("B1").Value = 2
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 1
$ws2.Range("A3").Value = 4
$ws2.Range("B3").Value = 3

# Conditional formatting: expression-based rules.
# Column A: highlight when A1<>1 (formula relative to the first cell
# of the range, like Excel does).
$rngA = $ws2.Range("A1:A3")
$fcA = $rngA.FormatConditions.Add(2, $null, "A1<>1")

# Column B: highlight when B1=1.
$rngB = $ws2.Range("B1:B3")
$fcB = $rngB.FormatConditions.Add(2, $null, "B1=1")

# Column B's rule was added second, but should end up with the higher
# priority (lower priority number) -- matches original authoring order
# captured in the fixture.
$rngB.FormatConditions.Item(1).SetFirstPriority() | Out-Null

# Selection / active sheet bookkeeping -- Sheet2 becomes the active tab.
$ws2.Range("C1").Select() | Out-Null
